$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '63.646.95'
$ws.Range("E2").Value = '  +0.32%  '

# Row 3
$ws.Range("D3").Value = '2.638.41'
$ws.Range("E3").Value = '  -1.40%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.56%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.23%  '

# Row 7
$ws.Range("E7").Value = '  +0.08%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.02%  '

# Row 9
$ws.Range("E9").Value = '  +1.68%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.382'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.53%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.60'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.16%  '

# Row 12
$ws.Range("E12").Value = '  -0.86%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.34'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.54%  '

# Row 14
$ws.Range("D14").Value = '3.114.48'
$ws.Range("E14").Value = '  -1.29%  '

# Row 15
$ws.Range("D15").Value = '63.488.75'
$ws.Range("E15").Value = '  +0.26%  '

# Row 16
$ws.Range("E16").Value = '  +0.97%  '

# Row 17
$ws.Range("D17").Value = '2.627.05'
$ws.Range("E17").Value = '  -1.96%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.74'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.38%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.56'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.18%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '345.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.15%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.87'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.35%  '

# Row 22
$ws.Range("E22").Value = '  -0.19%  '

# Row 23
$ws.Range("E23").Value = '  -2.11%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.95%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.66'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.78%  '

# Row 26
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.68'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.13%  '

# Row 27
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.18'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.74%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '559.33'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.08%  '

# Row 29
$ws.Range("E29").Value = '  +2.67%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.162'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.86%  '

# Row 31
$ws.Range("E31").Value = '  +0.14%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.05'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.46%  '

# Row 33
$ws.Range("D33").Value = '0.0₃0849'
$ws.Range("E33").Value = '  +4.88%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.75'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.75%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.26'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.20%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '169.62'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.71%  '

# Row 37
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.405'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.32%  '

# Row 38
$ws.Range("B38").Value = 'FirstDigitalUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.09%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.93'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.02%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.11'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.20%  '

# Row 41
$ws.Range("E41").Value = '  +0.08%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '164.77'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.45%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.04'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.38%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.79'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.94%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.80'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.89%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0565'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.06%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.628'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.17%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +14.06%  '

# Row 49
$ws.Range("E49").Value = '  +1.23%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0956'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.04%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.92%  '
